$wb = $excel.ActiveWorkbook

# Filtering of common results is now done by website title instead of urls:
# on the "DuckDuckGo" sheet, the keyword in A3 changes from "privacy" to
# "quality assurance" (so it matches the one already used on "Google").
$wsGoogle = $wb.Worksheets.Item("Google")
$wsDuck = $wb.Worksheets.Item("DuckDuckGo")

$wsDuck.Range("A3").Value = "quality assurance"

# Reproduce the resulting cursor/selection state left behind by the edit.
$wsGoogle.Activate()
[void]$wsGoogle.Range("A7").Select()

$wsDuck.Activate()
[void]$wsDuck.Range("A9").Select()
